# Weekly update: insert two new price observations into the Haba dataset.
#
# The sheet holds one row per (market, date) observation, each week a new
# observation is prepended to its block. This applies the "semanal"
# (weekly) update described in the commit message by inserting one new
# row above the old row 9 (pushing the rest of that block down by one)
# and a second new row above what was the old row 21 (pushing the tail
# of the block down by one more row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Insert-DataRow($RowIndex, $Fecha, $Volumen, $PrecioMinimo, $PrecioMaximo, $PrecioPromedio, $PrecioKg) {
    # Push the row currently at $RowIndex (and everything below) down by one.
    $ws.Rows.Item($RowIndex).Insert()

    $ws.Cells.Item($RowIndex, 1).Value = 2
    $ws.Cells.Item($RowIndex, 2).Value = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($RowIndex, 3).Value = "Coquimbo"

    $dateCell = $ws.Cells.Item($RowIndex, 4)
    $dateCell.Value = $Fecha
    $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($RowIndex, 5).Value = 4
    $ws.Cells.Item($RowIndex, 6).Value = 100112026
    $ws.Cells.Item($RowIndex, 7).Value = "Haba"
    $ws.Cells.Item($RowIndex, 8).Value = "Sin especificar"
    $ws.Cells.Item($RowIndex, 9).Value = "Primera"
    $ws.Cells.Item($RowIndex, 10).Value = $Volumen
    $ws.Cells.Item($RowIndex, 11).Value = $PrecioMinimo
    $ws.Cells.Item($RowIndex, 12).Value = $PrecioMaximo
    $ws.Cells.Item($RowIndex, 13).Value = $PrecioPromedio
    $ws.Cells.Item($RowIndex, 14).Value = "`$/saco 25 kilos"
    $ws.Cells.Item($RowIndex, 15).Value = "Provincia de Limarí"
    $ws.Cells.Item($RowIndex, 16).Value = $PrecioKg
    $ws.Cells.Item($RowIndex, 17).Value = 25
    $ws.Cells.Item($RowIndex, 18).Value = "Hortaliza"
}

# New observation dated 2021-08-26, inserted above the old row 9.
Insert-DataRow 9 44434 600 10000 11000 10500 420

# New observation dated 2021-08-27, inserted above what was row 21
# (now row 22, after the first insertion shifted it down by one).
Insert-DataRow 22 44435 600 10000 11000 10500 420
